# Insert two new weekly price rows at the top of the data block (rows 81-82),
# pushing all subsequent rows down by two. This matches the observed diff:
# dimension grows from R197 to R199, and every existing row r (>=83) ends up
# with the content that used to live at row r-2, with old rows 196/197
# reappearing unchanged as the new last rows 198/199.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before row 81 (shifts old row 81.. down to 83..)
$ws.Rows("81:82").Insert()

# New row 81: Coliflor "Primera" entry for 2021-10-14 (serial 44483)
$ws.Range("A81").Value = 4
$ws.Range("B81").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C81").Value = "Los Lagos"
$ws.Range("D81").Value = 44483
$ws.Range("E81").Value = 10
$ws.Range("F81").Value = 100112008
$ws.Range("G81").Value = "Coliflor"
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 250
$ws.Range("K81").Value = 1300
$ws.Range("L81").Value = 1300
$ws.Range("M81").Value = 1300
$ws.Range("N81").Value = "$/unidad"
$ws.Range("O81").Value = "Región Metropolitana"
$ws.Range("P81").Value = 1300
$ws.Range("Q81").Value = 1
$ws.Range("R81").Value = "Hortaliza"

# New row 82: Coliflor "Segunda" entry for the same date (2021-10-14)
$ws.Range("A82").Value = 4
$ws.Range("B82").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C82").Value = "Los Lagos"
$ws.Range("D82").Value = 44483
$ws.Range("E82").Value = 10
$ws.Range("F82").Value = 100112008
$ws.Range("G82").Value = "Coliflor"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Segunda"
$ws.Range("J82").Value = 250
$ws.Range("K82").Value = 1000
$ws.Range("L82").Value = 1000
$ws.Range("M82").Value = 1000
$ws.Range("N82").Value = "$/unidad"
$ws.Range("O82").Value = "Región Metropolitana"
$ws.Range("P82").Value = 1000
$ws.Range("Q82").Value = 1
$ws.Range("R82").Value = "Hortaliza"
